# Update gh-pages to output generated at 456a3b4
# Applies the numeric "want-to-go" count bumps and refreshed cover image
# links across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F7").Value = 98
$ws.Range("F8").Value = 10221
$ws.Range("F10").Value = 3525
$ws.Range("F12").Value = 2445
$ws.Range("F14").Value = 2807
$ws.Range("F17").Value = 2174
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202408/FBlScHDc1723775213878.jpeg"
$ws.Range("F21").Value = 386
$ws.Range("F23").Value = 149
$ws.Range("F24").Value = 315
$ws.Range("F26").Value = 224
$ws.Range("F30").Value = 1255
$ws.Range("F32").Value = 131
$ws.Range("F34").Value = 3707
$ws.Range("F35").Value = 3153
$ws.Range("F38").Value = 1044
$ws.Range("F42").Value = 98

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 16
$ws.Range("F8").Value = 7
$ws.Range("F16").Value = 179

# --- Sheet 3: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 748
$ws.Range("F3").Value = 988
$ws.Range("F4").Value = 126
$ws.Range("F5").Value = 2011

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 748
$ws.Range("F3").Value = 988
$ws.Range("F4").Value = 126
$ws.Range("F10").Value = 98
$ws.Range("F11").Value = 10221
$ws.Range("F13").Value = 3525
$ws.Range("F15").Value = 2445
$ws.Range("F17").Value = 16
$ws.Range("F19").Value = 2174
$ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202408/FBlScHDc1723775213878.jpeg"
$ws.Range("F23").Value = 386
$ws.Range("F24").Value = 149
$ws.Range("F25").Value = 315
$ws.Range("F26").Value = 224
$ws.Range("F29").Value = 1255
$ws.Range("F31").Value = 131
$ws.Range("F33").Value = 7
$ws.Range("F36").Value = 3153
$ws.Range("F37").Value = 1044
$ws.Range("F45").Value = 98
$ws.Range("F49").Value = 179
